{"js": "// The edited cells each hold a unique 'a\u00f7b=c, d' answer string in this\n// document, so we can safely locate and replace each one independently\n// with a scoped body.search() + insertText(\"Replace\") pass.\nconst replacements = [\n  [\"12\u00f74=3, 0\", \"20\u00f77=2, 6\"],\n  [\"59\u00f72=29, 1\", \"94\u00f76=15, 4\"],\n  [\"58\u00f76=9, 4\", \"35\u00f73=11, 2\"],\n  [\"51\u00f75=10, 1\", \"73\u00f79=8, 1\"],\n  [\"81\u00f75=16, 1\", \"92\u00f75=18, 2\"],\n  [\"52\u00f74=13, 0\", \"22\u00f74=5, 2\"],\n  [\"22\u00f78=2, 6\", \"92\u00f74=23, 0\"],\n  [\"34\u00f78=4, 2\", \"40\u00f74=10, 0\"],\n  [\"67\u00f72=33, 1\", \"92\u00f74=23, 0\"],\n  [\"24\u00f75=4, 4\", \"88\u00f73=29, 1\"],\n  [\"16\u00f72=8, 0\", \"21\u00f72=10, 1\"],\n  [\"55\u00f78=6, 7\", \"99\u00f78=12, 3\"],\n  [\"74\u00f73=24, 2\", \"49\u00f72=24, 1\"],\n  [\"26\u00f75=5, 1\", \"27\u00f77=3, 6\"],\n  [\"32\u00f74=8, 0\", \"54\u00f75=10, 4\"],\n  [\"59\u00f76=9, 5\", \"89\u00f74=22, 1\"],\n  [\"30\u00f78=3, 6\", \"25\u00f73=8, 1\"],\n  [\"99\u00f77=14, 1\", \"91\u00f77=13, 0\"],\n  [\"89\u00f77=12, 5\", \"64\u00f73=21, 1\"],\n  [\"18\u00f77=2, 4\", \"10\u00f79=1, 1\"],\n  [\"22\u00f79=2, 4\", \"95\u00f74=23, 3\"],\n  [\"56\u00f77=8, 0\", \"56\u00f72=28, 0\"],\n  [\"15\u00f79=1, 6\", \"20\u00f79=2, 2\"],\n  [\"50\u00f77=7, 1\", \"82\u00f72=41, 0\"],\n  [\"27\u00f75=5, 2\", \"65\u00f77=9, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Each old answer string is unique in the document, so a straightforward\n# Find/Replace-All pass per pair is safe and order-independent.\n$pairs = @(\n    @(\"12\u00f74=3, 0\", \"20\u00f77=2, 6\"),\n    @(\"59\u00f72=29, 1\", \"94\u00f76=15, 4\"),\n    @(\"58\u00f76=9, 4\", \"35\u00f73=11, 2\"),\n    @(\"51\u00f75=10, 1\", \"73\u00f79=8, 1\"),\n    @(\"81\u00f75=16, 1\", \"92\u00f75=18, 2\"),\n    @(\"52\u00f74=13, 0\", \"22\u00f74=5, 2\"),\n    @(\"22\u00f78=2, 6\", \"92\u00f74=23, 0\"),\n    @(\"34\u00f78=4, 2\", \"40\u00f74=10, 0\"),\n    @(\"67\u00f72=33, 1\", \"92\u00f74=23, 0\"),\n    @(\"24\u00f75=4, 4\", \"88\u00f73=29, 1\"),\n    @(\"16\u00f72=8, 0\", \"21\u00f72=10, 1\"),\n    @(\"55\u00f78=6, 7\", \"99\u00f78=12, 3\"),\n    @(\"74\u00f73=24, 2\", \"49\u00f72=24, 1\"),\n    @(\"26\u00f75=5, 1\", \"27\u00f77=3, 6\"),\n    @(\"32\u00f74=8, 0\", \"54\u00f75=10, 4\"),\n    @(\"59\u00f76=9, 5\", \"89\u00f74=22, 1\"),\n    @(\"30\u00f78=3, 6\", \"25\u00f73=8, 1\"),\n    @(\"99\u00f77=14, 1\", \"91\u00f77=13, 0\"),\n    @(\"89\u00f77=12, 5\", \"64\u00f73=21, 1\"),\n    @(\"18\u00f77=2, 4\", \"10\u00f79=1, 1\"),\n    @(\"22\u00f79=2, 4\", \"95\u00f74=23, 3\"),\n    @(\"56\u00f77=8, 0\", \"56\u00f72=28, 0\"),\n    @(\"15\u00f79=1, 6\", \"20\u00f79=2, 2\"),\n    @(\"50\u00f77=7, 1\", \"82\u00f72=41, 0\"),\n    @(\"27\u00f75=5, 2\", \"65\u00f77=9, 2\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $result) {\n        throw \"No match found for: $oldText\"\n    }\n}\n\n"}
